# Adds a new row to the "DBS" worksheet describing the
# findCustNoAndRecordDateFirst lookup function (key-read condition and
# order-by condition), matching the new SEQ entry added to the DB layout
# reference sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

$ws.Range("A8").Value = "findCustNoAndRecordDateFirst"
$ws.Range("B8").Value = "CustNo = ,AND RecordDate >= ,AND RecordDate <="
$ws.Range("C8").Value = "RecordDate ASC"

# Leave the selection where the author ended up after typing the new row.
[void]$ws.Range("C9").Select()
